# Update header labels on every worksheet in the workbook:
#   A1: "Input Sheet" -> "Car Name"
#   B1: "Value"       -> "Values"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Range("A1").Value2 -eq "Input Sheet") {
        $ws.Range("A1").Value = "Car Name"
    }
    if ($ws.Range("B1").Value2 -eq "Value") {
        $ws.Range("B1").Value = "Values"
    }
}
